$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 18.298577
$ws.Range("H2").Value = 54.895731
$ws.Range("I2").Value = 0.1985220285130613
$ws.Range("J2").Value = 0.1985220285130614
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 13.61877133333333
$ws.Range("N2").Value = 40.856314
$ws.Range("O2").Value = 0.3264056993691278
$ws.Range("P2").Value = 0.3264056993691277
$ws.Range("Q2").Value = 249.2041358883927
$ws.Range("R2").Value = 2242.837222995534
$ws.Range("S2").Value = 0.0647987215569837
$ws.Range("T2").Value = 0.06479872155698371

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 18.298577
$ws.Range("H3").Value = 54.895731
$ws.Range("I3").Value = 0.1985220285130613
$ws.Range("J3").Value = 0.1985220285130614
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 10.92748633333333
$ws.Range("N3").Value = 32.782459
$ws.Range("O3").Value = 0.2619027613928843
$ws.Range("P3").Value = 0.2619027613928842
$ws.Range("Q3").Value = 199.9574500869477
$ws.Range("R3").Value = 1799.617050782529
$ws.Range("S3").Value = 0.05199346746488766
$ws.Range("T3").Value = 0.05199346746488766

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 18.298577
$ws.Range("H4").Value = 54.895731
$ws.Range("I4").Value = 0.1985220285130613
$ws.Range("J4").Value = 0.1985220285130614
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 7.492675333333334
$ws.Range("N4").Value = 22.478026
$ws.Range("O4").Value = 0.1795794842620271
$ws.Range("P4").Value = 0.1795794842620271
$ws.Range("Q4").Value = 137.1052965230006
$ws.Range("R4").Value = 1233.947668707006
$ws.Range("S4").Value = 0.03565048349502699
$ws.Range("T4").Value = 0.03565048349502699

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 18.298577
$ws.Range("H5").Value = 54.895731
$ws.Range("I5").Value = 0.1985220285130613
$ws.Range("J5").Value = 0.1985220285130614
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 9.684515333333332
$ws.Range("N5").Value = 29.053546
$ws.Range("O5").Value = 0.2321120549759609
$ws.Range("P5").Value = 0.2321120549759609
$ws.Range("Q5").Value = 177.2128495346806
$ws.Range("R5").Value = 1594.915645812126
$ws.Range("S5").Value = 0.04607935599616297
$ws.Range("T5").Value = 0.04607935599616297

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 12.03986033333333
$ws.Range("H6").Value = 36.119581
$ws.Range("I6").Value = 0.1306209491802164
$ws.Range("J6").Value = 0.1306209491802164
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 13.61877133333333
$ws.Range("N6").Value = 40.856314
$ws.Range("O6").Value = 0.3264056993691278
$ws.Range("P6").Value = 0.3264056993691277
$ws.Range("Q6").Value = 163.9681047649371
$ws.Range("R6").Value = 1475.712942884434
$ws.Range("S6").Value = 0.04263542226942782
$ws.Range("T6").Value = 0.04263542226942782

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 12.03986033333333
$ws.Range("H7").Value = 36.119581
$ws.Range("I7").Value = 0.1306209491802164
$ws.Range("J7").Value = 0.1306209491802164
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 10.92748633333333
$ws.Range("N7").Value = 32.782459
$ws.Range("O7").Value = 0.2619027613928843
$ws.Range("P7").Value = 0.2619027613928842
$ws.Range("Q7").Value = 131.5654092477421
$ws.Range("R7").Value = 1184.088683229679
$ws.Range("S7").Value = 0.03420998728605826
$ws.Range("T7").Value = 0.03420998728605826

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 12.03986033333333
$ws.Range("H8").Value = 36.119581
$ws.Range("I8").Value = 0.1306209491802164
$ws.Range("J8").Value = 0.1306209491802164
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 7.492675333333334
$ws.Range("N8").Value = 22.478026
$ws.Range("O8").Value = 0.1795794842620271
$ws.Range("P8").Value = 0.1795794842620271
$ws.Range("Q8").Value = 90.2107645363451
$ws.Range("R8").Value = 811.896880827106
$ws.Range("S8").Value = 0.0234568426875997
$ws.Range("T8").Value = 0.0234568426875997

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 12.03986033333333
$ws.Range("H9").Value = 36.119581
$ws.Range("I9").Value = 0.1306209491802164
$ws.Range("J9").Value = 0.1306209491802164
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 9.684515333333332
$ws.Range("N9").Value = 29.053546
$ws.Range("O9").Value = 0.2321120549759609
$ws.Range("P9").Value = 0.2321120549759609
$ws.Range("Q9").Value = 116.6002120093584
$ws.Range("R9").Value = 1049.401908084226
$ws.Range("S9").Value = 0.03031869693713058
$ws.Range("T9").Value = 0.03031869693713058

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 54.94736600000001
$ws.Range("H10").Value = 164.842098
$ws.Range("I10").Value = 0.5961262758174921
$ws.Range("J10").Value = 0.5961262758174922
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 13.61877133333333
$ws.Range("N10").Value = 40.856314
$ws.Range("O10").Value = 0.3264056993691278
$ws.Range("P10").Value = 0.3264056993691277
$ws.Range("Q10").Value = 748.3156129229749
$ws.Range("R10").Value = 6734.840516306774
$ws.Range("S10").Value = 0.1945790139705221
$ws.Range("T10").Value = 0.1945790139705221

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 54.94736600000001
$ws.Range("H11").Value = 164.842098
$ws.Range("I11").Value = 0.5961262758174921
$ws.Range("J11").Value = 0.5961262758174922
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 10.92748633333333
$ws.Range("N11").Value = 32.782459
$ws.Range("O11").Value = 0.2619027613928843
$ws.Range("P11").Value = 0.2619027613928842
$ws.Range("Q11").Value = 600.4365910176648
$ws.Range("R11").Value = 5403.929319158983
$ws.Range("S11").Value = 0.1561271177754573
$ws.Range("T11").Value = 0.1561271177754573

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 54.94736600000001
$ws.Range("H12").Value = 164.842098
$ws.Range("I12").Value = 0.5961262758174921
$ws.Range("J12").Value = 0.5961262758174922
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 7.492675333333334
$ws.Range("N12").Value = 22.478026
$ws.Range("O12").Value = 0.1795794842620271
$ws.Range("P12").Value = 0.1795794842620271
$ws.Range("Q12").Value = 411.7027738598388
$ws.Range("R12").Value = 3705.324964738548
$ws.Range("S12").Value = 0.1070520491663482
$ws.Range("T12").Value = 0.1070520491663481

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 54.94736600000001
$ws.Range("H13").Value = 164.842098
$ws.Range("I13").Value = 0.5961262758174921
$ws.Range("J13").Value = 0.5961262758174922
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 9.684515333333332
$ws.Range("N13").Value = 29.053546
$ws.Range("O13").Value = 0.2321120549759609
$ws.Range("P13").Value = 0.2321120549759609
$ws.Range("Q13").Value = 532.1386085532787
$ws.Range("R13").Value = 4789.247476979508
$ws.Range("S13").Value = 0.1383680949051646
$ws.Range("T13").Value = 0.1383680949051646

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 6.888234666666667
$ws.Range("H14").Value = 20.664704
$ws.Range("I14").Value = 0.07473074648923014
$ws.Range("J14").Value = 0.07473074648923016
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 13.61877133333333
$ws.Range("N14").Value = 40.856314
$ws.Range("O14").Value = 0.3264056993691278
$ws.Range("P14").Value = 0.3264056993691277
$ws.Range("Q14").Value = 93.8092928156729
$ws.Range("R14").Value = 844.2836353410561
$ws.Range("S14").Value = 0.02439254157219416
$ws.Range("T14").Value = 0.02439254157219416

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 6.888234666666667
$ws.Range("H15").Value = 20.664704
$ws.Range("I15").Value = 0.07473074648923014
$ws.Range("J15").Value = 0.07473074648923016
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 10.92748633333333
$ws.Range("N15").Value = 32.782459
$ws.Range("O15").Value = 0.2619027613928843
$ws.Range("P15").Value = 0.2619027613928842
$ws.Range("Q15").Value = 75.2710901807929
$ws.Range("R15").Value = 677.4398116271361
$ws.Range("S15").Value = 0.01957218886648096
$ws.Range("T15").Value = 0.01957218886648096

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 6.888234666666667
$ws.Range("H16").Value = 20.664704
$ws.Range("I16").Value = 0.07473074648923014
$ws.Range("J16").Value = 0.07473074648923016
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 7.492675333333334
$ws.Range("N16").Value = 22.478026
$ws.Range("O16").Value = 0.1795794842620271
$ws.Range("P16").Value = 0.1795794842620271
$ws.Range("Q16").Value = 51.61130597714489
$ws.Range("R16").Value = 464.501753794304
$ws.Range("S16").Value = 0.01342010891305224
$ws.Range("T16").Value = 0.01342010891305224

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 6.888234666666667
$ws.Range("H17").Value = 20.664704
$ws.Range("I17").Value = 0.07473074648923014
$ws.Range("J17").Value = 0.07473074648923016
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 9.684515333333332
$ws.Range("N17").Value = 29.053546
$ws.Range("O17").Value = 0.2321120549759609
$ws.Range("P17").Value = 0.2321120549759609
$ws.Range("Q17").Value = 66.70921424893154
$ws.Range("R17").Value = 600.382928240384
$ws.Range("S17").Value = 0.01734590713750279
$ws.Range("T17").Value = 0.01734590713750279

